$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.520102666666667
$ws.Range("H2").Value = 13.560308
$ws.Range("I2").Value = 0.9927775608668273
$ws.Range("J2").Value = 0.9927775608668273
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08378199999999998
$ws.Range("N2").Value = 0.251346
$ws.Range("O2").Value = 0.007571394704126512
$ws.Range("P2").Value = 0.007571394704126512
$ws.Range("Q2").Value = 0.3787032416186666
$ws.Range("R2").Value = 3.408329174568
$ws.Range("S2").Value = 0.007516710766722732
$ws.Range("T2").Value = 0.007516710766722732

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.520102666666667
$ws.Range("H3").Value = 13.560308
$ws.Range("I3").Value = 0.9927775608668273
$ws.Range("J3").Value = 0.9927775608668273
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.07352966666666667
$ws.Range("N3").Value = 0.220589
$ws.Range("O3").Value = 0.006644889460697858
$ws.Range("P3").Value = 0.006644889460697857
$ws.Range("Q3").Value = 0.3323616423791112
$ws.Range("R3").Value = 2.991254781412001
$ws.Range("S3").Value = 0.006596897151021307
$ws.Range("T3").Value = 0.006596897151021306

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.520102666666667
$ws.Range("H4").Value = 13.560308
$ws.Range("I4").Value = 0.9927775608668273
$ws.Range("J4").Value = 0.9927775608668273
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.908285
$ws.Range("N4").Value = 32.724855
$ws.Range("O4").Value = 0.9857837158351757
$ws.Range("P4").Value = 0.9857837158351755
$ws.Range("Q4").Value = 49.30656811726001
$ws.Range("R4").Value = 443.75911305534
$ws.Range("S4").Value = 0.9786639529490834
$ws.Range("T4").Value = 0.9786639529490833

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03288366666666667
$ws.Range("H5").Value = 0.098651
$ws.Range("I5").Value = 0.007222439133172593
$ws.Range("J5").Value = 0.007222439133172593
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08378199999999998
$ws.Range("N5").Value = 0.251346
$ws.Range("O5").Value = 0.007571394704126512
$ws.Range("P5").Value = 0.007571394704126512
$ws.Range("Q5").Value = 0.002755059360666666
$ws.Range("R5").Value = 0.024795534246
$ws.Range("S5").Value = 0.00005468393740377904
$ws.Range("T5").Value = 0.00005468393740377904

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.03288366666666667
$ws.Range("H6").Value = 0.098651
$ws.Range("I6").Value = 0.007222439133172593
$ws.Range("J6").Value = 0.007222439133172593
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.07352966666666667
$ws.Range("N6").Value = 0.220589
$ws.Range("O6").Value = 0.006644889460697858
$ws.Range("P6").Value = 0.006644889460697857
$ws.Range("Q6").Value = 0.002417925048777778
$ws.Range("R6").Value = 0.021761325439
$ws.Range("S6").Value = 0.00004799230967655033
$ws.Range("T6").Value = 0.00004799230967655032

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03288366666666667
$ws.Range("H7").Value = 0.098651
$ws.Range("I7").Value = 0.007222439133172593
$ws.Range("J7").Value = 0.007222439133172593
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.908285
$ws.Range("N7").Value = 32.724855
$ws.Range("O7").Value = 0.9857837158351757
$ws.Range("P7").Value = 0.9857837158351755
$ws.Range("Q7").Value = 0.358704407845
$ws.Range("R7").Value = 3.228339670605
$ws.Range("S7").Value = 0.007119762886092263
$ws.Range("T7").Value = 0.007119762886092263
